$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1046.3857
$ws.Range("J17").Value = 1054.7537
$ws.Range("L17").Value = 3164.2611
$ws.Range("N17").Value = -3500.2611
$ws.Range("H74").Value = 16531.889
$ws.Range("I74").Value = 16531.889
$ws.Range("K74").Value = 16531.889
$ws.Range("M74").Value = -15595.889
$ws.Range("H77").Value = 16531.889
$ws.Range("I77").Value = 16531.889
$ws.Range("K77").Value = 82659.44499999999
$ws.Range("M77").Value = -77979.44499999999
$ws.Range("H87").Value = 75000
$ws.Range("J87").Value = 75000
$ws.Range("L87").Value = 75000
$ws.Range("N87").Value = -77496
$ws.Range("H90").Value = 75000
$ws.Range("J90").Value = 75000
$ws.Range("L90").Value = 225000
$ws.Range("N90").Value = -237480
$ws.Range("H116").Value = 5884.143
$ws.Range("I116").Value = 5531.5
$ws.Range("K116").Value = 5531.5
$ws.Range("M116").Value = -2089.5
$ws.Range("H118").Value = 1425.7646
$ws.Range("I118").Value = 1153.8
$ws.Range("J118").Value = 1814.2858
$ws.Range("K118").Value = 3461.4
$ws.Range("L118").Value = 5442.857400000001
$ws.Range("M118").Value = -1804.4
$ws.Range("N118").Value = -8756.857400000001
$ws.Range("H137").Value = 1470.1852
$ws.Range("I137").Value = 1479.8334
$ws.Range("J137").Value = 1393
$ws.Range("K137").Value = 4439.5002
$ws.Range("L137").Value = 4179
$ws.Range("M137").Value = -1889.5002
$ws.Range("N137").Value = -9279
$ws.Range("H138").Value = 1680.125
$ws.Range("I138").Value = 1240.7894
$ws.Range("J138").Value = 3349.6
$ws.Range("K138").Value = 3722.3682
$ws.Range("L138").Value = 10048.8
$ws.Range("M138").Value = 1417.6318
$ws.Range("N138").Value = -20328.8
$ws.Range("H141").Value = 3112
$ws.Range("I141").Value = 2470.6155
$ws.Range("K141").Value = 7411.8465
$ws.Range("M141").Value = -2231.8465

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 29365.455
$ws.Range("I37").Value = 21666.666
$ws.Range("J37").Value = 32252.5
$ws.Range("K37").Value = 21666.666
$ws.Range("L37").Value = 32252.5
$ws.Range("M37").Value = -21393.666
$ws.Range("N37").Value = -32798.5
$ws.Range("H45").Value = 1673.75
$ws.Range("I45").Value = 1673.75
$ws.Range("K45").Value = 1673.75
$ws.Range("M45").Value = -1296.75
$ws.Range("H50").Value = 1922
$ws.Range("I50").Value = 694.5
$ws.Range("J50").Value = 3763.25
$ws.Range("K50").Value = 694.5
$ws.Range("L50").Value = 3763.25
$ws.Range("M50").Value = 19.5
$ws.Range("N50").Value = -5191.25

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 13000
$ws.Range("I33").Value = 13000
$ws.Range("K33").Value = 13000
$ws.Range("M33").Value = -12664
$ws.Range("H86").Value = 33337466
$ws.Range("I86").Value = 3998.7
$ws.Range("K86").Value = 3998.7
$ws.Range("M86").Value = -2875.7
$ws.Range("H89").Value = 33337466
$ws.Range("I89").Value = 3998.7
$ws.Range("K89").Value = 19993.5
$ws.Range("M89").Value = -14377.5
$ws.Range("H105").Value = 3191.577
$ws.Range("I105").Value = 2547.0952
$ws.Range("K105").Value = 2547.0952
$ws.Range("M105").Value = -800.0952000000002

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1895.8572
$ws.Range("J22").Value = 2134.2
$ws.Range("L22").Value = 2134.2
$ws.Range("N22").Value = -2834.2
$ws.Range("H35").Value = 299
$ws.Range("I35").Value = 299
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 299
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -5
$ws.Range("N35").ClearContents()
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H132").Value = 2050.1738
$ws.Range("I132").Value = 1715.4117
$ws.Range("K132").Value = 5146.2351
$ws.Range("M132").Value = -2616.2351

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 5900.7646
$ws.Range("J23").Value = 4094.8
$ws.Range("L23").Value = 12284.4
$ws.Range("N23").Value = -12754.4
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H131").Value = 2173.3635
$ws.Range("I131").Value = 1474
$ws.Range("K131").Value = 4422
$ws.Range("M131").Value = 618

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 7500
$ws.Range("J27").Value = 7500
$ws.Range("L27").Value = 7500
$ws.Range("N27").Value = -7832
$ws.Range("H70").Value = 9340.5
$ws.Range("J70").Value = 10004.5
$ws.Range("L70").Value = 10004.5
$ws.Range("N70").Value = -10544.5
$ws.Range("H73").Value = 9340.5
$ws.Range("J73").Value = 10004.5
$ws.Range("L73").Value = 10004.5
$ws.Range("N73").Value = -11876.5
$ws.Range("H102").Value = 16499.666
$ws.Range("I102").Value = 9500
$ws.Range("J102").Value = 19999.5
$ws.Range("K102").Value = 9500
$ws.Range("L102").Value = 19999.5
$ws.Range("M102").Value = -7878
$ws.Range("N102").Value = -23243.5
$ws.Range("H113").Value = 7299.8
$ws.Range("I113").Value = 3908.6365
$ws.Range("K113").Value = 3908.6365
$ws.Range("M113").Value = -1738.6365
$ws.Range("H122").Value = 22460.105
$ws.Range("I122").Value = 23916.268
$ws.Range("K122").Value = 71748.804
$ws.Range("M122").Value = -69298.804
$ws.Range("H126").Value = 3050
$ws.Range("I126").Value = 3100
$ws.Range("K126").Value = 9300
$ws.Range("M126").Value = -6830
$ws.Range("H132").Value = 2191.7144
$ws.Range("I132").Value = 2170.8823
$ws.Range("J132").Value = 2900
$ws.Range("K132").Value = 6512.646900000001
$ws.Range("L132").Value = 8700
$ws.Range("M132").Value = -3982.646900000001
$ws.Range("N132").Value = -13760

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2519.6667
$ws.Range("I16").Value = 295
$ws.Range("J16").Value = 6969
$ws.Range("K16").Value = 295
$ws.Range("L16").Value = 6969
$ws.Range("M16").Value = -125
$ws.Range("N16").Value = -7309
$ws.Range("H136").Value = 12030.472
$ws.Range("I136").Value = 1713.9
$ws.Range("K136").Value = 5141.700000000001
$ws.Range("M136").Value = -2591.700000000001
$ws.Range("H139").Value = 99998
$ws.Range("J139").Value = 99998
$ws.Range("L139").Value = 99998
$ws.Range("N139").Value = -110278

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 47500
$ws.Range("I99").Value = 45000
$ws.Range("K99").Value = 45000
$ws.Range("M99").Value = -42005
$ws.Range("H126").Value = 2049.2727
$ws.Range("I126").Value = 1441.8572
$ws.Range("K126").Value = 4325.571599999999
$ws.Range("M126").Value = -1855.571599999999
$ws.Range("H132").Value = 1377.5454
$ws.Range("I132").Value = 1395.0476
$ws.Range("J132").Value = 1346.9166
$ws.Range("K132").Value = 4185.142800000001
$ws.Range("L132").Value = 4040.7498
$ws.Range("M132").Value = -1655.142800000001
$ws.Range("N132").Value = -9100.7498
